$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: date, calidad, volumen, precios, unidad, precio/kg, kg/unidad
$ws.Range("D2").Value = 44495
$ws.Range("L2").Value = "Segunda"
$ws.Range("M2").Value = 270
$ws.Range("N2").Value = 19000
$ws.Range("O2").Value = 20000
$ws.Range("P2").Value = 19556
$ws.Range("Q2").Value = "`$/bandeja 8 kilos"
$ws.Range("S2").Value = 2444
$ws.Range("T2").Value = 8

# Row 4: date, volumen, precios, precio/kg
$ws.Range("D4").Value = 44488
$ws.Range("M4").Value = 160
$ws.Range("N4").Value = 17000
$ws.Range("O4").Value = 18000
$ws.Range("P4").Value = 17500
$ws.Range("S4").Value = 2188

# Row 5: date, calidad, volumen, precios, unidad, precio/kg, kg/unidad
$ws.Range("D5").Value = 44162
$ws.Range("L5").Value = "Primera"
$ws.Range("M5").Value = 200
$ws.Range("N5").Value = 2000
$ws.Range("O5").Value = 2100
$ws.Range("P5").Value = 2050
$ws.Range("Q5").Value = "`$/kilo (en caja de 14 kilos)"
$ws.Range("S5").Value = 2050
$ws.Range("T5").Value = 1
